$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: shift each row's data up by one (row N <- old row N+1),
# dropping the oldest entry (old row 31) and appending a new weekly entry at row 41.

# Row 31
$ws.Range("D31").Value = 44245
$ws.Range("H31").Value = 'Patagonia'
$ws.Range("I31").Value = '1a (cosecha)'
$ws.Range("K31").Value = 8500
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 8750
$ws.Range("N31").Value = '$/saco 25 kilos'
$ws.Range("O31").Value = 'Región de Los Lagos'
$ws.Range("P31").Value = 350

# Row 32
$ws.Range("D32").Value = 44350
$ws.Range("H32").Value = 'Asterix'
$ws.Range("I32").Value = '1a (cosecha lavada)'
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 9500
$ws.Range("N32").Value = '$/malla 25 kilos'
$ws.Range("O32").Value = 'Región de La Araucanía'
$ws.Range("P32").Value = 380

# Row 33
$ws.Range("D33").Value = 44574
$ws.Range("I33").Value = '1a (cosecha)'
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 14000
$ws.Range("M33").Value = 13500
$ws.Range("O33").Value = 'Región del Maule'
$ws.Range("P33").Value = 540

# Row 34
$ws.Range("D34").Value = 44601
$ws.Range("I34").Value = '1a (cosecha lavada)'
$ws.Range("K34").Value = 11000
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = 11500
$ws.Range("N34").Value = '$/saco 25 kilos'
$ws.Range("O34").Value = 'Región de Los Lagos'
$ws.Range("P34").Value = 460

# Row 35
$ws.Range("H35").Value = 'Patagonia'
$ws.Range("I35").Value = '1a (cosecha)'
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9500
$ws.Range("O35").Value = 'Provincia de Melipilla'
$ws.Range("P35").Value = 380

# Row 36
$ws.Range("D36").Value = 44592
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = 10500
$ws.Range("O36").Value = 'Región del Maule'
$ws.Range("P36").Value = 420

# Row 37
$ws.Range("D37").Value = 44433
$ws.Range("H37").Value = 'Asterix'
$ws.Range("I37").Value = '1a (guarda)'
$ws.Range("K37").Value = 9000
$ws.Range("L37").Value = 9500
$ws.Range("M37").Value = 9250
$ws.Range("O37").Value = 'Región de La Araucanía'
$ws.Range("P37").Value = 370

# Row 38
$ws.Range("D38").Value = 44159
$ws.Range("H38").Value = 'Rosara'
$ws.Range("I38").Value = '1a (cosecha lavada)'
$ws.Range("K38").Value = 11000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 11500
$ws.Range("N38").Value = '$/malla 25 kilos'
$ws.Range("O38").Value = 'Región del Maule'
$ws.Range("P38").Value = 460

# Row 39
$ws.Range("D39").Value = 44386
$ws.Range("H39").Value = 'Asterix'
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 9500
$ws.Range("M39").Value = 9250
$ws.Range("O39").Value = 'Región de Los Lagos'
$ws.Range("P39").Value = 370

# Row 40
$ws.Range("D40").Value = 44322
$ws.Range("H40").Value = 'Rodeo'
$ws.Range("I40").Value = '1a (cosecha)'
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = 8250
$ws.Range("N40").Value = '$/saco 25 kilos'
$ws.Range("P40").Value = 330

# Row 41
$ws.Range("D41").Value = 44651
$ws.Range("H41").Value = 'Asterix'
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 9500
$ws.Range("P41").Value = 380
